$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("相談件数")

# The sheet currently ends with a footer row at 118 (A118 empty, B118 holds the
# "※4/8…" note). A new data row for 5/22 (serial 43973) needs to be inserted
# just above that footer, pushing the footer down to row 119.

# 1) Push the footer row (A118:B118) down to row 119, keeping its formatting.
$ws.Range("A118:B118").Copy($ws.Range("A119:B119"))

# 2) Give the new row 118 the same number formats as the row above it (117)
#    before writing values, so dates/numbers render the same way.
$ws.Range("A117:E117").Copy($ws.Range("A118:E118"))

# 3) Write the new day's figures into row 118.
$ws.Range("A118").Value = 43973
$ws.Range("B118").Value = 136
$ws.Range("C118").Value = 38601
$ws.Range("D118").Value = 37
$ws.Range("E118").Value = 7801

# 4) Extend the print area to cover the new last row.
foreach ($n in $wb.Names) {
    if ($n.Name.Contains("Print_Area")) {
        $n.RefersTo = "=相談件数!`$A`$1:`$E`$119"
    }
}

# 5) Move the active selection to where the footer label now sits, matching
#    where the editor's cursor ended up after the insert.
$ws.Activate()
$ws.Range("B119").Select()
